$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(20, 1).Value = "CB187"
$ws.Cells.Item(20, 2).Value = 0.66
$ws.Cells.Item(20, 3).Value = 0.84
$ws.Cells.Item(20, 4).Value = 0.9399999999999999

$ws.Cells.Item(21, 1).Value = "CB194"
$ws.Cells.Item(21, 2).Value = 0.5
$ws.Cells.Item(21, 3).Value = 0.66
$ws.Cells.Item(21, 4).Value = 1.14
